$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: Several "Price" cells below hold numeric-looking text (e.g. "0.999",
# "7.30") that must stay as literal strings, matching the source workbook's
# inline-string cells. A leading apostrophe forces Excel to store the text
# literally instead of auto-converting it to a number; the apostrophe itself
# is not stored as part of the cell's value/text.

$ws.Range("D2").Value = "57.007.47"
$ws.Range("E2").Value = "  -1.55%  "
$ws.Range("D3").Value = "2.990.78"
$ws.Range("E3").Value = "  -1.91%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'500.32"
$ws.Range("E5").Value = "  -4.50%  "
$ws.Range("D6").Value = "'138.35"
$ws.Range("E6").Value = "  -2.68%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("E8").Value = "  -3.32%  "
$ws.Range("D9").Value = "'7.30"
$ws.Range("E9").Value = "  -4.11%  "
$ws.Range("E10").Value = "  -3.69%  "
$ws.Range("D11").Value = "'0.360"
$ws.Range("E11").Value = "  -2.33%  "
$ws.Range("D12").Value = "3.496.00"
$ws.Range("E12").Value = "  -2.15%  "
$ws.Range("E13").Value = "  -2.30%  "
$ws.Range("D14").Value = "'26.20"
$ws.Range("E14").Value = "  -2.67%  "
$ws.Range("D15").Value = "'0.0000161"
$ws.Range("E15").Value = "  -6.40%  "
$ws.Range("D16").Value = "57.054.94"
$ws.Range("E16").Value = "  -1.41%  "
$ws.Range("D17").Value = "'6.11"
$ws.Range("E17").Value = "  -1.82%  "
$ws.Range("D18").Value = "2.990.70"
$ws.Range("E18").Value = "  -2.14%  "
$ws.Range("D19").Value = "'12.69"
$ws.Range("E19").Value = "  -1.86%  "
$ws.Range("D20").Value = "'7.91"
$ws.Range("E20").Value = "  -3.60%  "
$ws.Range("D21").Value = "'323.16"
$ws.Range("E21").Value = "  -5.04%  "
$ws.Range("E22").Value = "  -0.18%  "
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("D24").Value = "'0.493"
$ws.Range("E24").Value = "  -1.20%  "
$ws.Range("D25").Value = "'63.51"
$ws.Range("E25").Value = "  -2.48%  "
$ws.Range("E26").Value = "  +0.32%  "
$ws.Range("E27").Value = "  -4.82%  "
$ws.Range("D28").Value = "0.0₃0897"
$ws.Range("E28").Value = "  -7.59%  "
$ws.Range("D29").Value = "'6.58"
$ws.Range("E29").Value = "  -5.97%  "
$ws.Range("D30").Value = "'7.11"
$ws.Range("E30").Value = "  -2.48%  "
$ws.Range("D31").Value = "'1.77"
$ws.Range("E31").Value = "  -4.60%  "
$ws.Range("D32").Value = "'1.17"
$ws.Range("E32").Value = "  -6.26%  "
$ws.Range("D33").Value = "'20.34"
$ws.Range("E33").Value = "  -3.67%  "
$ws.Range("D34").Value = "'155.58"
$ws.Range("E34").Value = "  -0.57%  "
$ws.Range("D35").Value = "'4.59"
$ws.Range("E35").Value = "  -2.82%  "
$ws.Range("D36").Value = "'5.82"
$ws.Range("E36").Value = "  -1.20%  "
$ws.Range("E37").Value = "  -5.95%  "
$ws.Range("D38").Value = "'24.37"
$ws.Range("E38").Value = "  -5.82%  "
$ws.Range("D39").Value = "'0.0668"
$ws.Range("E39").Value = "  -4.10%  "
$ws.Range("B40").Value = "OKB"
$ws.Range("C40").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D40").Value = "'37.83"
$ws.Range("E40").Value = "  +0.34%  "
$ws.Range("B41").Value = "RenzoRestakedETH"
$ws.Range("C41").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D41").Value = "3.015.86"
$ws.Range("E41").Value = "  -2.22%  "
$ws.Range("E42").Value = "  -0.09%  "
$ws.Range("D43").Value = "'3.77"
$ws.Range("E43").Value = "  -2.57%  "
$ws.Range("D44").Value = "'0.646"
$ws.Range("E44").Value = "  -2.62%  "
$ws.Range("D45").Value = "2.208.09"
$ws.Range("E45").Value = "  -5.28%  "
$ws.Range("E46").Value = "  -5.63%  "
$ws.Range("B47").Value = "ONDO"
$ws.Range("C47").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D47").Value = "'0.948"
$ws.Range("E47").Value = "  -8.37%  "
$ws.Range("B48").Value = "Cosmos"
$ws.Range("C48").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D48").Value = "'5.99"
$ws.Range("E48").Value = "  -0.74%  "
$ws.Range("D49").Value = "'0.0235"
$ws.Range("E49").Value = "  -5.15%  "
$ws.Range("D50").Value = "'19.37"
$ws.Range("E50").Value = "  -3.31%  "
$ws.Range("E51").Value = "  -10.92%  "
